# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Jengibre" (Vega Modelo de Temuco) right
# after the existing row 321, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 322 (pushes old row 322.. down to 323..)
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(322, 1).Value2 = 10
$ws.Cells.Item(322, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(322, 3).Value2 = "La Araucanía"
$ws.Cells.Item(322, 4).Value2 = 45275
$ws.Cells.Item(322, 5).Value2 = 9
$ws.Cells.Item(322, 6).Value2 = 100114007
$ws.Cells.Item(322, 7).Value2 = "Jengibre"
$ws.Cells.Item(322, 8).Value2 = "Sin especificar"
$ws.Cells.Item(322, 9).Value2 = "Primera"
$ws.Cells.Item(322, 10).Value2 = 35
$ws.Cells.Item(322, 11).Value2 = 25000
$ws.Cells.Item(322, 12).Value2 = 25000
$ws.Cells.Item(322, 13).Value2 = 25000
$ws.Cells.Item(322, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(322, 15).Value2 = "Perú"
$ws.Cells.Item(322, 16).Value2 = 1923
$ws.Cells.Item(322, 17).Value2 = 13
$ws.Cells.Item(322, 18).Value2 = "Hortaliza"
